$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 115.27273
$ws.Range("I6").Value = 116.7
$ws.Range("K6").Value = 350.1
$ws.Range("M6").Value = -238.1

$ws.Range("H15").Value = 1181.3889
$ws.Range("I15").Value = 1181.3889
$ws.Range("K15").Value = 3544.1667
$ws.Range("M15").Value = -3375.1667

$ws.Range("H80").Value = 1086.2916
$ws.Range("I80").Value = 1139.5294
$ws.Range("K80").Value = 3418.5882
$ws.Range("M80").Value = -2420.5882

$ws.Range("H83").Value = 1086.2916
$ws.Range("I83").Value = 1139.5294
$ws.Range("K83").Value = 10255.7646
$ws.Range("M83").Value = -5263.764599999999

$ws.Range("H106").Value = 32175.646
$ws.Range("I106").Value = 34461.23
$ws.Range("J106").Value = 24747.5
$ws.Range("K106").Value = 34461.23
$ws.Range("L106").Value = 24747.5
$ws.Range("M106").Value = -33830.23
$ws.Range("N106").Value = -26009.5

$ws.Range("H116").Value = 6373.5
$ws.Range("I116").Value = 4747.5
$ws.Range("J116").Value = 7999.5
$ws.Range("K116").Value = 4747.5
$ws.Range("L116").Value = 7999.5
$ws.Range("M116").Value = -1305.5
$ws.Range("N116").Value = -14883.5

$ws.Range("H138").Value = 3893.2622
$ws.Range("J138").Value = 4622.5713
$ws.Range("L138").Value = 13867.7139
$ws.Range("N138").Value = -24147.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4499.8335
$ws.Range("I45").Value = 4499.8335
$ws.Range("K45").Value = 4499.8335
$ws.Range("M45").Value = -4122.8335

$ws.Range("H132").Value = 874.1111
$ws.Range("I132").Value = 858.375
$ws.Range("K132").Value = 2575.125
$ws.Range("M132").Value = -45.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1928.5555
$ws.Range("I20").Value = 1787.7142
$ws.Range("K20").Value = 1787.7142
$ws.Range("M20").Value = -1540.7142

$ws.Range("H75").Value = 44999.332
$ws.Range("I75").Value = 15000
$ws.Range("K75").Value = 15000
$ws.Range("M75").Value = -14064

$ws.Range("H78").Value = 44999.332
$ws.Range("I78").Value = 15000
$ws.Range("K78").Value = 45000
$ws.Range("M78").Value = -40320

$ws.Range("H88").Value = 34556
$ws.Range("J88").Value = 34556
$ws.Range("L88").Value = 34556
$ws.Range("N88").Value = -35368

$ws.Range("H91").Value = 34556
$ws.Range("J91").Value = 34556
$ws.Range("L91").Value = 34556
$ws.Range("N91").Value = -37364

$ws.Range("H94").Value = 1657.3529
$ws.Range("J94").Value = 2480
$ws.Range("L94").Value = 2480
$ws.Range("N94").Value = -3382

$ws.Range("H99").Value = 2570.2942
$ws.Range("I99").Value = 2174.25
$ws.Range("K99").Value = 2174.25
$ws.Range("M99").Value = -676.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2821.6
$ws.Range("I2").Value = 1025.75
$ws.Range("K2").Value = 1025.75
$ws.Range("M2").Value = -912.75

$ws.Range("H16").Value = 3816
$ws.Range("J16").Value = 3749
$ws.Range("L16").Value = 3749
$ws.Range("N16").Value = -4323

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H88").Value = 95430.336
$ws.Range("J88").Value = 95430.336
$ws.Range("L88").Value = 95430.336
$ws.Range("N88").Value = -96242.336

$ws.Range("H91").Value = 95430.336
$ws.Range("J91").Value = 95430.336
$ws.Range("L91").Value = 95430.336
$ws.Range("N91").Value = -98238.336

$ws.Range("H98").Value = 125000
$ws.Range("J98").Value = 125000
$ws.Range("L98").Value = 125000
$ws.Range("N98").Value = -129492

$ws.Range("H113").Value = 3816
$ws.Range("J113").Value = 3749
$ws.Range("L113").Value = 3749
$ws.Range("N113").Value = -8089

$ws.Range("H122").Value = 3356.3845
$ws.Range("I122").Value = 3577.2727
$ws.Range("J122").Value = 2141.5
$ws.Range("K122").Value = 10731.8181
$ws.Range("L122").Value = 6424.5
$ws.Range("M122").Value = -8281.8181
$ws.Range("N122").Value = -11324.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 3538.1875
$ws.Range("J62").Value = 3038.3845
$ws.Range("L62").Value = 9115.1535
$ws.Range("N62").Value = -10487.1535

$ws.Range("H65").Value = 3538.1875
$ws.Range("J65").Value = 3038.3845
$ws.Range("L65").Value = 27345.4605
$ws.Range("N65").Value = -34209.4605

$ws.Range("H97").Value = 8000
$ws.Range("I97").Value = 8000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 24000
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -23504
$ws.Range("N97").ClearContents()

$ws.Range("H99").Value = 660
$ws.Range("I99").Value = 660
$ws.Range("K99").Value = 1980
$ws.Range("M99").Value = 266

$ws.Range("H114").Value = 200
$ws.Range("I114").Value = 200
$ws.Range("J114").Value = 200
$ws.Range("K114").Value = 600
$ws.Range("L114").Value = 600
$ws.Range("M114").Value = 2654
$ws.Range("N114").Value = -7108

$ws.Range("H128").Value = 670773.5
$ws.Range("I128").Value = 670773.5
$ws.Range("K128").Value = 2012320.5
$ws.Range("M128").Value = -2007340.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 2333.3333
$ws.Range("I22").Value = 2250
$ws.Range("K22").Value = 2250
$ws.Range("M22").Value = -1721

$ws.Range("H97").Value = 444.33334
$ws.Range("I97").Value = 375
$ws.Range("K97").Value = 375
$ws.Range("M97").Value = 121

$ws.Range("H99").Value = 10442.75
$ws.Range("I99").Value = 1808.4
$ws.Range("J99").Value = 24833.334
$ws.Range("K99").Value = 1808.4
$ws.Range("L99").Value = 24833.334
$ws.Range("M99").Value = 437.5999999999999
$ws.Range("N99").Value = -29325.334

$ws.Range("H102").Value = 2574.5715
$ws.Range("I102").Value = 2520.3333
$ws.Range("K102").Value = 2520.3333
$ws.Range("M102").Value = -898.3332999999998

$ws.Range("H113").Value = 2281.75
$ws.Range("I113").Value = 1625.6666
$ws.Range("J113").Value = 4250
$ws.Range("K113").Value = 1625.6666
$ws.Range("L113").Value = 4250
$ws.Range("M113").Value = 544.3334
$ws.Range("N113").Value = -8590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4827.4287
$ws.Range("J7").Value = 4449
$ws.Range("L7").Value = 4449
$ws.Range("N7").Value = -4673

$ws.Range("H22").Value = 3806
$ws.Range("I22").Value = 2739.4
$ws.Range("K22").Value = 2739.4
$ws.Range("M22").Value = -2444.4

$ws.Range("H27").Value = 3806
$ws.Range("I27").Value = 2739.4
$ws.Range("K27").Value = 2739.4
$ws.Range("M27").Value = -2632.4

$ws.Range("H40").Value = 2705.4666
$ws.Range("I40").Value = 2814.5386
$ws.Range("J40").Value = 1996.5
$ws.Range("K40").Value = 2814.5386
$ws.Range("L40").Value = 1996.5
$ws.Range("M40").Value = -2678.5386
$ws.Range("N40").Value = -2268.5

$ws.Range("H46").Value = 5141.5
$ws.Range("I46").Value = 1998.75
$ws.Range("K46").Value = 1998.75
$ws.Range("M46").Value = -1810.75

$ws.Range("H61").Value = 3727.1
$ws.Range("I61").Value = 3708.375
$ws.Range("K61").Value = 3708.375
$ws.Range("M61").Value = -3506.375

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H113").Value = 3727.1
$ws.Range("I113").Value = 3708.375
$ws.Range("K113").Value = 3708.375
$ws.Range("M113").Value = -1538.375

$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 4000
$ws.Range("K122").Value = 12000
$ws.Range("M122").Value = -9550

$ws.Range("H126").Value = 4827.4287
$ws.Range("J126").Value = 4449
$ws.Range("L126").Value = 13347
$ws.Range("N126").Value = -18287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 17483.334
$ws.Range("J69").Value = 17483.334
$ws.Range("L69").Value = 17483.334
$ws.Range("N69").Value = -18981.334

$ws.Range("H72").Value = 17483.334
$ws.Range("J72").Value = 17483.334
$ws.Range("L72").Value = 52450.00199999999
$ws.Range("N72").Value = -59938.00199999999

$ws.Range("H113").Value = 815.8333
$ws.Range("I113").Value = 833
$ws.Range("J113").Value = 798.6667
$ws.Range("K113").Value = 2499
$ws.Range("L113").Value = 2396.0001
$ws.Range("M113").Value = -329
$ws.Range("N113").Value = -6736.0001

$ws.Range("H122").Value = 3166.6667
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
